$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2371541501976284
$ws.Range("C2").Value = 0.466403162055336
$ws.Range("J2").Value = 0.007905138339920948
$ws.Range("P2").Value = 0.1857707509881423
$ws.Range("S2").Value = 0.1027667984189723
$ws.Range("B3").Value = 0.008264462809917356
$ws.Range("C3").Value = 0.02479338842975207
$ws.Range("J3").Value = 0.02479338842975207
$ws.Range("P3").Value = 0.8264462809917356
$ws.Range("S3").Value = 0.115702479338843
$ws.Range("P4").Value = 0.5769230769230769
$ws.Range("S4").Value = 0.4230769230769231
$ws.Range("B6").Value = 0.05392156862745098
$ws.Range("D6").Value = 0.009803921568627451
$ws.Range("F6").Value = 0.07352941176470588
$ws.Range("J6").Value = 0.2352941176470588
$ws.Range("O6").Value = 0.02450980392156863
$ws.Range("Q6").Value = 0.1372549019607843
$ws.Range("R6").Value = 0.107843137254902
$ws.Range("S6").Value = 0.357843137254902
$ws.Range("B7").Value = 0.1395348837209302
$ws.Range("D7").Value = 0.005813953488372093
$ws.Range("E7").Value = 0.01162790697674419
$ws.Range("F7").Value = 0.06395348837209303
$ws.Range("J7").Value = 0.1395348837209302
$ws.Range("O7").Value = 0.01162790697674419
$ws.Range("Q7").Value = 0.1395348837209302
$ws.Range("R7").Value = 0.08139534883720931
$ws.Range("S7").Value = 0.4069767441860465
$ws.Range("B8").Value = 0.078125
$ws.Range("D8").Value = 0.01339285714285714
$ws.Range("E8").Value = 0.002232142857142857
$ws.Range("F8").Value = 0.05803571428571429
$ws.Range("J8").Value = 0.09375
$ws.Range("O8").Value = 0.006696428571428571
$ws.Range("Q8").Value = 0.1941964285714286
$ws.Range("R8").Value = 0.078125
$ws.Range("S8").Value = 0.4754464285714285
$ws.Range("B9").Value = 0.09795918367346938
$ws.Range("D9").Value = 0.01224489795918367
$ws.Range("E9").Value = 0.004081632653061225
$ws.Range("F9").Value = 0.04081632653061224
$ws.Range("J9").Value = 0.1142857142857143
$ws.Range("O9").Value = 0.00816326530612245
$ws.Range("Q9").Value = 0.1591836734693877
$ws.Range("R9").Value = 0.06938775510204082
$ws.Range("S9").Value = 0.4938775510204081
$ws.Range("B10").Value = 0.07762938230383973
$ws.Range("D10").Value = 0.01001669449081803
$ws.Range("E10").Value = 0.0008347245409015025
$ws.Range("F10").Value = 0.07595993322203673
$ws.Range("J10").Value = 0.09766277128547579
$ws.Range("O10").Value = 0.01335559265442404
$ws.Range("Q10").Value = 0.2220367278797997
$ws.Range("R10").Value = 0.07262103505843072
$ws.Range("S10").Value = 0.4298831385642738
$ws.Range("G11").Value = 0.1515151515151515
$ws.Range("J11").Value = 0.07575757575757576
$ws.Range("K11").Value = 0.196969696969697
$ws.Range("L11").Value = 0.5681818181818182
$ws.Range("S11").Value = 0.007575757575757576
$ws.Range("G12").Value = 0.7115384615384616
$ws.Range("J12").Value = 0.1858974358974359
$ws.Range("K12").Value = 0.00641025641025641
$ws.Range("L12").Value = 0.04487179487179487
$ws.Range("S12").Value = 0.05128205128205128
$ws.Range("G13").Value = 0.625
$ws.Range("J13").Value = 0.325
$ws.Range("S13").Value = 0.05
$ws.Range("F15").Value = 0.01818181818181818
$ws.Range("H15").Value = 0.1636363636363636
$ws.Range("I15").Value = 0.1136363636363636
$ws.Range("J15").Value = 0.3727272727272727
$ws.Range("K15").Value = 0.05
$ws.Range("M15").Value = 0.01818181818181818
$ws.Range("N15").Value = 0.004545454545454545
$ws.Range("O15").Value = 0.05
$ws.Range("S15").Value = 0.2090909090909091
$ws.Range("H16").Value = 0.15527950310559
$ws.Range("I16").Value = 0.1055900621118012
$ws.Range("J16").Value = 0.5031055900621118
$ws.Range("K16").Value = 0.08695652173913043
$ws.Range("M16").Value = 0.0124223602484472
$ws.Range("O16").Value = 0.04347826086956522
$ws.Range("S16").Value = 0.09316770186335403
$ws.Range("F17").Value = 0.01360544217687075
$ws.Range("H17").Value = 0.1836734693877551
$ws.Range("I17").Value = 0.09977324263038549
$ws.Range("J17").Value = 0.4081632653061225
$ws.Range("K17").Value = 0.08390022675736962
$ws.Range("M17").Value = 0.02267573696145125
$ws.Range("N17").Value = 0.00453514739229025
$ws.Range("O17").Value = 0.08390022675736962
$ws.Range("S17").Value = 0.09977324263038549
$ws.Range("F18").Value = 0.02298850574712644
$ws.Range("H18").Value = 0.2241379310344828
$ws.Range("I18").Value = 0.1149425287356322
$ws.Range("J18").Value = 0.396551724137931
$ws.Range("K18").Value = 0.07471264367816093
$ws.Range("M18").Value = 0.01149425287356322
$ws.Range("O18").Value = 0.09770114942528736
$ws.Range("S18").Value = 0.05747126436781609
$ws.Range("F19").Value = 0.01489028213166144
$ws.Range("H19").Value = 0.2100313479623825
$ws.Range("I19").Value = 0.109717868338558
$ws.Range("J19").Value = 0.3816614420062696
$ws.Range("K19").Value = 0.103448275862069
$ws.Range("M19").Value = 0.01724137931034483
$ws.Range("N19").Value = 0.002351097178683386
$ws.Range("O19").Value = 0.06739811912225706
$ws.Range("S19").Value = 0.0932601880877743
